# Atualização automática de GARIBALDI.xlsx
#
# - Renomeia a planilha "Paineis DARQ" para "PAINEIS DARQ"
# - Renomeia a planilha "Recolhimento x Eliminacao" para "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove a planilha "Desarquivamentos Pendentes" (não é mais utilizada)

$wb = $excel.ActiveWorkbook

# Avoid any confirmation prompt when deleting a worksheet
$excel.DisplayAlerts = $false | Out-Null

# Rename sheets
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the obsolete sheet entirely
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null

$wb.Save()
